$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "E.g.: lets say ... 4 ways/blocks." -> "... 4 ways or blocks."
#    Replace the whole sentence (this also clears the proofErr gramStart/
#    gramEnd marks that wrapped "lets" and merges the sentence into a single
#    run, matching how Word re-flows a run after an in-place edit).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "E.g.: lets say we have 1 set & a 4-way associative cache which means each set has 4 ways/blocks.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "E.g.: lets say we have 1 set & a 4-way associative cache which means each set has 4 ways or blocks.",
    2) | Out-Null

# Re-split the merged run into "...4 ways" | " " | "or " | "blocks." and drop
# the (zero width) _GoBack bookmark exactly between the space and "or " --
# this mirrors where Word leaves _GoBack after the last edit.
$rng = $d.Content
$rng.Find.Execute("or blocks.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$orStart = $rng.Start
$blocksStart = $orStart + 3
$spaceStart = $orStart - 1

$bm = $d.Range($blocksStart, $blocksStart)
$bm.Bookmarks.Add("TempSplit1") | Out-Null
$d.Bookmarks("TempSplit1").Delete()

$bm = $d.Range($spaceStart, $spaceStart)
$bm.Bookmarks.Add("TempSplit2") | Out-Null
$d.Bookmarks("TempSplit2").Delete()

# Move (or create) the _GoBack bookmark to sit right before "or " -- Word only
# keeps a single _GoBack bookmark, so this also removes it from wherever it
# used to be (originally before " 2-bit counter" in the "d)" paragraph).
$bm = $d.Range($orStart, $orStart)
$bm.Bookmarks.Add("_GoBack") | Out-Null

# ---------------------------------------------------------------------------
# 2) "b)Way" + "1" + "(Block" + "1" + ") requires " + "one" + " 2-bit counter"
#    -> "b)Way" (kept wrapped in the original gramStart/gramEnd proofErr) +
#    "1(Block1) requires one 2-bit counter" (merged) + " [Set0]" (untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("b)Way", $true, $false, $false, $false, $false, $true, 1, $false, "b)Way", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("1(Block1) requires one 2-bit counter [Set0]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$setStart = $rng.End - 7

$bm = $d.Range($setStart, $setStart)
$bm.Bookmarks.Add("TempSplitB") | Out-Null
$d.Bookmarks("TempSplitB").Delete()

# ---------------------------------------------------------------------------
# 3) "c" + ")Way" + "2" + "(Block" + "2" + ") requires " + "one" + " 2-bit counter"
#    -> "c)Way2(Block2) requires one 2-bit counter" (merged, no proofErr) +
#    " [Set 0]" (kept separate).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("c)Way2(Block2) requires one 2-bit counter", $true, $false, $false, $false, $false, $true, 1, $false, "c)Way2(Block2) requires one 2-bit counter", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("c)Way2(Block2) requires one 2-bit counter [Set 0]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$setStart = $rng.End - 8

$bm = $d.Range($setStart, $setStart)
$bm.Bookmarks.Add("TempSplitC") | Out-Null
$d.Bookmarks("TempSplitC").Delete()

# ---------------------------------------------------------------------------
# 4) "d)" + "Way" + "3" + "(Block" + "3" + ") requires " + "one" + " 2-bit counter"
#    -> "d)Way3(Block3) requires " (merged) + "one" (kept separate) +
#    " 2-bit counter" (kept separate, _GoBack already relocated above) +
#    " [Set 0]" (untouched).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("d)Way3(Block3) requires ", $true, $false, $false, $false, $false, $true, 1, $false, "d)Way3(Block3) requires ", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("d)Way3(Block3) requires one 2-bit counter", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$reqEnd = $rng.Start + 24
$oneEnd = $rng.Start + 27

$bm = $d.Range($oneEnd, $oneEnd)
$bm.Bookmarks.Add("TempSplitD1") | Out-Null
$d.Bookmarks("TempSplitD1").Delete()

$bm = $d.Range($reqEnd, $reqEnd)
$bm.Bookmarks.Add("TempSplitD2") | Out-Null
$d.Bookmarks("TempSplitD2").Delete()
